# Automatische test-sync: 2025-06-22 18:44:50
# Adds the new "Productinformatie" mail-log entry (row 14) to the Logs sheet,
# adds the matching summary row (row 11) to the Dashboard sheet, and extends
# the chart series + conditional-formatting ranges to cover the new rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Logs sheet: append row 14
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A14").Value = "Productinformatie"
$logs.Range("B14").Value = "mailmind.test@zohomail.eu"
$logs.Range("C14").Value = "Wat is het verschil tussen product A en product B?"
$logs.Range("D14").Value = "Productinformatie"
$logs.Range("E14").Value = "Beste klant,`nDank u voor uw interesse in onze producten A en B. Het belangrijkste verschil tussen product A en B is de functionaliteit. Product A is ontworpen met functies zoals X, Y en Z, terwijl product B functies heeft zoals P, Q en R. Afhankelijk van uw behoeften en vereisten kunt u het product kiezen dat het beste bij uw situatie past.`nAls u meer gedetailleerde informatie nodig heeft over specifieke functies of kenmerken van product A en B, dan helpen wij u graag verder.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F14").Value = "2025-06-22 18:44:13"
$logs.Range("G14").Value = "Ja"

# Extend the two conditional-formatting blocks (D and G columns) so they
# cover the newly added row 14 as well.
$logs.Range("D2:D13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D14"))
$logs.Range("G2:G13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G14"))

# ---------------------------------------------------------------------------
# 2. Dashboard sheet: append summary row 11
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A11").Value = "Productinformatie"
$dash.Range("B11").Value = 1

# ---------------------------------------------------------------------------
# 3. Chart: extend category/value series ranges to include the new row
# ---------------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$11,Dashboard!`$B`$2:`$B`$11,1)"
